# Update the "Förändrad" (changed) date column for all existing data rows
# (rows 2-235) from 2023-10-03 (45202) to 2023-10-04 (45203).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C235").Value = 45203

# Row 235 picks up an explicit row height in the new file.
$ws.Rows.Item(235).RowHeight = 15

# Append the new record (row 236) for case "A 47393-2023", copying the
# formatting (date number format on B/C, wrap-text style on R) from the
# row above it, then filling in the values.
$ws.Range("A235:R235").Copy()
$ws.Range("A236").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A236").Value = "A 47393-2023"
$ws.Range("B236").Value = 45202
$ws.Range("C236").Value = 45203
$ws.Range("D236").Value = "HALLANDS LÄN"
$ws.Range("E236").Value = "KUNGSBACKA"
$ws.Range("F236").Clear()
$ws.Range("G236").Value = 0.7
$ws.Range("H236").Value = 0
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = 0
$ws.Range("K236").Value = 0
$ws.Range("L236").Value = 0
$ws.Range("M236").Value = 0
$ws.Range("N236").Value = 0
$ws.Range("O236").Value = 0
$ws.Range("P236").Value = 0
$ws.Range("Q236").Value = 0
$ws.Range("R236").Value = ""
